$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 96

# Copy formatting from the row above (row 95) so style indexes match (bold/border on A, date format on E)
$ws.Range("A95:V95").Copy() | Out-Null
$ws.Range("A" + $newRow + ":V" + $newRow).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($newRow, 1).Value = 95
$ws.Cells.Item($newRow, 2).Value = "ecuador"
$ws.Cells.Item($newRow, 3).Value = "liga-pro"

# "2023" must land as a text cell (like the other text columns), but a
# plain string Value assignment gets auto-detected as the number 2023, and
# forcing text via NumberFormat/quote-prefix leaves a permanent, unused
# style behind. Route it through a TEXT() formula, then paste-special the
# result back as a value so the cell keeps its default (unstyled) text type.
$ws.Cells.Item($newRow, 4).Formula = "=TEXT(2023,""0"")"
$ws.Cells.Item($newRow, 4).Copy() | Out-Null
$ws.Cells.Item($newRow, 4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item($newRow, 5).Value = 45237.04166666666
$ws.Cells.Item($newRow, 6).Value = "U. Catolica"
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = "Libertad"
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 1.44
$ws.Cells.Item($newRow, 11).Value = "02/11/2023 01:12"
$ws.Cells.Item($newRow, 12).Value = 1.37
$ws.Cells.Item($newRow, 13).Value = "07/11/2023 00:56"
$ws.Cells.Item($newRow, 14).Value = 4.45
$ws.Cells.Item($newRow, 15).Value = "02/11/2023 01:12"
$ws.Cells.Item($newRow, 16).Value = 5.07
$ws.Cells.Item($newRow, 17).Value = "07/11/2023 00:56"
$ws.Cells.Item($newRow, 18).Value = 6.11
$ws.Cells.Item($newRow, 19).Value = "02/11/2023 01:12"
$ws.Cells.Item($newRow, 20).Value = 8.07
$ws.Cells.Item($newRow, 21).Value = "07/11/2023 00:56"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/ecuador/liga-pro/u-catolica-libertad/Cz8gzlzh/"
